$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab/name from "SA-HW50.xpc" to "SA"
$ws.Name = "SA"

# Add a new data row (row 16) mirroring the scheme rows above it:
# HKL-style index column, label column, and 14 columns of "1"s.

# Copy formatting (bold, border, centered) from the row above (A15) so the
# new index cell A16 matches the existing style used for column A entries.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
